# Auto-generated edit script applying numeric updates from the diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 142.73334
$ws.Range("I33").Value = 142.73334
$ws.Range("K33").Value = 142.73334
$ws.Range("M33").Value = 86.26666
$ws.Range("H112").Value = 3105.7144
$ws.Range("J112").Value = 3198.5186
$ws.Range("L112").Value = 9595.5558
$ws.Range("N112").Value = -11811.5558
$ws.Range("H116").Value = 3204.5652
$ws.Range("I116").Value = 1519.4
$ws.Range("K116").Value = 1519.4
$ws.Range("M116").Value = 1922.6
$ws.Range("H129").Value = 833917.25
$ws.Range("J129").Value = 1429330
$ws.Range("L129").Value = 4287990
$ws.Range("N129").Value = -4297990

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1441.4736
$ws.Range("I2").Value = 1423.6666
$ws.Range("K2").Value = 1423.6666
$ws.Range("M2").Value = -1310.6666
$ws.Range("H32").Value = 3674.1738
$ws.Range("I32").Value = 2978.2195
$ws.Range("J32").Value = 9381
$ws.Range("K32").Value = 2978.2195
$ws.Range("L32").Value = 9381
$ws.Range("M32").Value = -2691.2195
$ws.Range("N32").Value = -9955
$ws.Range("H45").Value = 2595.3333
$ws.Range("I45").Value = 3587.5557
$ws.Range("K45").Value = 3587.5557
$ws.Range("M45").Value = -3210.5557
$ws.Range("H74").Value = 76924216
$ws.Range("I74").Value = 125000670
$ws.Range("J74").Value = 1890
$ws.Range("K74").Value = 125000670
$ws.Range("L74").Value = 1890
$ws.Range("M74").Value = -124999796
$ws.Range("N74").Value = -3638
$ws.Range("H77").Value = 76924216
$ws.Range("I77").Value = 125000670
$ws.Range("J77").Value = 1890
$ws.Range("K77").Value = 625003350
$ws.Range("L77").Value = 9450
$ws.Range("M77").Value = -624998982
$ws.Range("N77").Value = -18186
$ws.Range("H116").Value = 1441.4736
$ws.Range("I116").Value = 1423.6666
$ws.Range("K116").Value = 1423.6666
$ws.Range("M116").Value = 870.3334
$ws.Range("H122").Value = 2318.1428
$ws.Range("I122").Value = 1887.9166
$ws.Range("K122").Value = 5663.7498
$ws.Range("M122").Value = -3213.7498
$ws.Range("H132").Value = 15073.868
$ws.Range("I132").Value = 1949.8077
$ws.Range("J132").Value = 43509.332
$ws.Range("K132").Value = 5849.4231
$ws.Range("L132").Value = 130527.996
$ws.Range("M132").Value = -3319.4231
$ws.Range("N132").Value = -135587.996

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1441.4736
$ws.Range("I3").Value = 1423.6666
$ws.Range("K3").Value = 1423.6666
$ws.Range("M3").Value = -1309.6666
$ws.Range("H94").Value = 2541.15
$ws.Range("I94").Value = 2394.5625
$ws.Range("K94").Value = 2394.5625
$ws.Range("M94").Value = -1943.5625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 851
$ws.Range("I16").Value = 808.3
$ws.Range("K16").Value = 808.3
$ws.Range("M16").Value = -521.3
$ws.Range("H58").Value = 23104.61
$ws.Range("I58").Value = 1533
$ws.Range("J58").Value = 63551.375
$ws.Range("K58").Value = 1533
$ws.Range("L58").Value = 63551.375
$ws.Range("M58").Value = -1330
$ws.Range("N58").Value = -63957.375
$ws.Range("H113").Value = 851
$ws.Range("I113").Value = 808.3
$ws.Range("K113").Value = 808.3
$ws.Range("M113").Value = 1361.7
$ws.Range("H136").Value = 23104.61
$ws.Range("I136").Value = 1533
$ws.Range("J136").Value = 63551.375
$ws.Range("K136").Value = 4599
$ws.Range("L136").Value = 190654.125
$ws.Range("M136").Value = -2049
$ws.Range("N136").Value = -195754.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1083
$ws.Range("I68").Value = 750
$ws.Range("J68").Value = 1249.5
$ws.Range("K68").Value = 2250
$ws.Range("L68").Value = 3748.5
$ws.Range("M68").Value = -1439
$ws.Range("N68").Value = -5370.5
$ws.Range("H71").Value = 1083
$ws.Range("I71").Value = 750
$ws.Range("J71").Value = 1249.5
$ws.Range("K71").Value = 6750
$ws.Range("L71").Value = 11245.5
$ws.Range("M71").Value = -2694
$ws.Range("N71").Value = -19357.5
$ws.Range("H122").Value = 938
$ws.Range("J122").Value = 1048.6
$ws.Range("L122").Value = 9437.4
$ws.Range("N122").Value = -14337.4
$ws.Range("H131").Value = 751.49
$ws.Range("J131").Value = 751.49
$ws.Range("L131").Value = 2254.47
$ws.Range("N131").Value = -12334.47

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 2986.6667
$ws.Range("I53").Value = 1980
$ws.Range("J53").Value = 5000
$ws.Range("K53").Value = 1980
$ws.Range("L53").Value = 5000
$ws.Range("M53").Value = -1349
$ws.Range("N53").Value = -6262
$ws.Range("H113").Value = 2918.5625
$ws.Range("I113").Value = 2429.7
$ws.Range("K113").Value = 2429.7
$ws.Range("M113").Value = -259.6999999999998
$ws.Range("H122").Value = 63493310
$ws.Range("I122").Value = 19609136
$ws.Range("K122").Value = 58827408
$ws.Range("M122").Value = -58824958
$ws.Range("H131").Value = 37663
$ws.Range("J131").Value = 37663
$ws.Range("L131").Value = 37663
$ws.Range("N131").Value = -47743
$ws.Range("H132").Value = 21788.822
$ws.Range("I132").Value = 5114.2
$ws.Range("J132").Value = 41028.77
$ws.Range("K132").Value = 15342.6
$ws.Range("L132").Value = 123086.31
$ws.Range("M132").Value = -12812.6
$ws.Range("N132").Value = -128146.31

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1228145.2
$ws.Range("I122").Value = 1510703.4
$ws.Range("J122").Value = 3726.6667
$ws.Range("K122").Value = 4532110.199999999
$ws.Range("L122").Value = 11180.0001
$ws.Range("M122").Value = -4529660.199999999
$ws.Range("N122").Value = -16080.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H82").Value = 39250.5
$ws.Range("J82").Value = 39250.5
$ws.Range("L82").Value = 39250.5
$ws.Range("N82").Value = -40016.5
$ws.Range("H85").Value = 39250.5
$ws.Range("J85").Value = 39250.5
$ws.Range("L85").Value = 39250.5
$ws.Range("N85").Value = -41902.5
$ws.Range("H122").Value = 1857.579
$ws.Range("I122").Value = 1896.875
$ws.Range("J122").Value = 1648
$ws.Range("K122").Value = 5690.625
$ws.Range("L122").Value = 4944
$ws.Range("M122").Value = -3240.625
$ws.Range("N122").Value = -9844
